$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg)
# for rows 2-18. Values are a weekly re-shuffle of the same data set.

$data = @{
    2  = @{ D = 44895; M = 240; N = 3000;  O = 3500;  P = 3250;  S = 1625 }
    3  = @{ D = 44517; M = 400; N = 5500;  O = 6000;  P = 5750;  S = 2875 }
    5  = @{ D = 44455; M = 200; N = 12000; O = 13000; P = 12500; S = 6250 }
    6  = @{ D = 44475; M = 240; N = 11000; O = 12000; P = 11500; S = 5750 }
    7  = @{ D = 44490; M = 400; N = 9500;  O = 10000; P = 9750;  S = 4875 }
    8  = @{ D = 44875; M = 400; N = 7000;  O = 7500;  P = 7250;  S = 3625 }
    9  = @{ D = 44889; M = 460; N = 3500;  O = 4000;  P = 3750;  S = 1875 }
    10 = @{ D = 44881; M = 440; N = 6000;  O = 7000;  P = 6500;  S = 3250 }
    12 = @{ D = 44461; M = 200; N = 11000; O = 12000; P = 11500; S = 5750 }
    13 = @{ D = 44882; M = 440; N = 6000;  O = 7000;  P = 6500;  S = 3250 }
    14 = @{ D = 44482; M = 240; N = 10000; O = 11000; P = 10500; S = 5250 }
    15 = @{ D = 44497; M = 500; N = 9000;  O = 10000; P = 9500;  S = 4750 }
    16 = @{ D = 44454; M = 160; N = 12000; O = 13000; P = 12500; S = 6250 }
    17 = @{ D = 44489; M = 160; N = 9500;  O = 10000; P = 9750;  S = 4875 }
    18 = @{ D = 44819; M = 240; N = 11000; O = 12000; P = 11500; S = 5750 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("S$row").Value = $vals.S
}
